{"js": "const body = context.document.body;\nconst paras = body.paragraphs;\nparas.load(\"text,style\");\nawait context.sync();\n\n// Locate the \"Overview\" Heading 1 paragraph and remove it, along with the\n// now-redundant blank paragraph that preceded it.\nlet overviewIdx = -1;\nfor (let i = 0; i < paras.items.length; i++) {\n  if (paras.items[i].text === \"Overview\") {\n    overviewIdx = i;\n    break;\n  }\n}\nif (overviewIdx !== -1) {\n  paras.items[overviewIdx].delete();\n  if (overviewIdx - 1 >= 0 && paras.items[overviewIdx - 1].text === \"\") {\n    paras.items[overviewIdx - 1].delete();\n  }\n  await context.sync();\n}\n\n// Change \"Sprint #1\" heading into two runs: \"Sprint \" + \"1\" (drop the \"#\").\nparas.load(\"text,style\");\nawait context.sync();\nlet sprintIdx = -1;\nfor (let i = 0; i < paras.items.length; i++) {\n  if (paras.items[i].text === \"Sprint #1\") {\n    sprintIdx = i;\n    break;\n  }\n}\nif (sprintIdx !== -1) {\n  const target = paras.items[sprintIdx];\n  const ooxml =\n    '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    \"<pkg:xmlData>\" +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    \"<w:body>\" +\n    '<w:p><w:pPr><w:pStyle w:val=\"Heading1\"/></w:pPr>' +\n    '<w:r><w:t xml:space=\"preserve\">Sprint </w:t></w:r>' +\n    \"<w:r><w:t>1</w:t></w:r>\" +\n    \"</w:p>\" +\n    \"</w:body></w:document>\" +\n    \"</pkg:xmlData></pkg:part></pkg:package>\";\n  target.insertOoxml(ooxml, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n", "ps1": "$d = $word.ActiveDocument\n\n# --- Remove the \"Overview\" Heading 1 paragraph, along with the blank\n#     paragraph that used to separate it from the title above it. ---\n$overviewIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Text.TrimEnd(\"`r\") -eq \"Overview\") {\n        $overviewIndex = $i\n        break\n    }\n}\nif ($overviewIndex -gt 0) {\n    $d.Paragraphs.Item($overviewIndex).Range.Delete()\n    $prevIndex = $overviewIndex - 1\n    if ($prevIndex -ge 1) {\n        $prevPara = $d.Paragraphs.Item($prevIndex)\n        if ($prevPara.Range.Text.TrimEnd(\"`r\") -eq \"\") {\n            $prevPara.Range.Delete()\n        }\n    }\n}\n\n# --- \"Sprint #1\" -> two runs: \"Sprint \" and \"1\" (drop the \"#\"). ---\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$found = $rng.Find.Execute(\"Sprint #1\")\nif ($found) {\n    $rng.Text = \"\"\n    $xml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n      '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n      '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n      '<pkg:xmlData>' +\n      '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n      '<w:body>' +\n      '<w:p><w:pPr><w:pStyle w:val=\"Heading1\"/></w:pPr>' +\n      '<w:r><w:t xml:space=\"preserve\">Sprint </w:t></w:r>' +\n      '<w:r><w:t>1</w:t></w:r>' +\n      '</w:p>' +\n      '</w:body></w:document>' +\n      '</pkg:xmlData></pkg:part></pkg:package>'\n    $rng.InsertXML($xml)\n}\n"}
